$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.596.62'
$ws.Range('E2').Value = '  +0.08%  '
$ws.Range('D3').Value = '1.889.28'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '244.85'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.95%  '
$ws.Range('E6').Value = '  +0.03%  '
$ws.Range('E7').Value = '  -0.43%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2915'
$ws.Range('D8').ClearFormats()
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06500'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.36%  '
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07764'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.66%  '
$ws.Range('D12').Value = '1.888.37'
$ws.Range('E12').Value = '  +0.86%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7384'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -0.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '96.18'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.37%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.196'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '284.44'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('D17').Value = '30.660.53'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.06'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('E19').Value = '  +0.08%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000007505'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = '2.135.56'
$ws.Range('E21').Value = '  +0.91%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.002'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.20%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.271'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.264'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +1.17%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.176'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.68%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '164.34'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.87'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.907'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.68%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.349'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09760'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -2.31%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.477'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.90%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.300'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +0.05%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.138'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +0.65%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04883'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.129'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.92%  '
$ws.Range('E36').Value = '  -0.54%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.709'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.28%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01901'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +2.67%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.838'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +3.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '75.56'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +2.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.206'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +0.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.010'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +2.05%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4267'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +2.01%  '
$ws.Range('E44').Value = '  +0.05%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8257'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.94%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '101.54'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.532'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.60%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.44'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +0.16%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.982'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.26%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '910.73'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05752'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +1.85%  '
